$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-10
# from the old date serial (45224 -> 2023-10-25) to the new one
# (45233 -> 2023-11-03).
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
